$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-25 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-26 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("62÷9=6, 8", $true, $false, $false, $false, $false, $true, 1, $false, "61÷6=10, 1", 2) | Out-Null
$d.Content.Find.Execute("59÷3=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "35÷3=11, 2", 2) | Out-Null
$d.Content.Find.Execute("64÷2=32, 0", $true, $false, $false, $false, $false, $true, 1, $false, "58÷2=29, 0", 2) | Out-Null
$d.Content.Find.Execute("89÷7=12, 5", $true, $false, $false, $false, $false, $true, 1, $false, "48÷3=16, 0", 2) | Out-Null
$d.Content.Find.Execute("51÷7=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "23÷4=5, 3", 2) | Out-Null
$d.Content.Find.Execute("30÷8=3, 6", $true, $false, $false, $false, $false, $true, 1, $false, "14÷7=2, 0", 2) | Out-Null
$d.Content.Find.Execute("48÷4=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "79÷8=9, 7", 2) | Out-Null
$d.Content.Find.Execute("33÷4=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "22÷9=2, 4", 2) | Out-Null
$d.Content.Find.Execute("20÷6=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "19÷7=2, 5", 2) | Out-Null
$d.Content.Find.Execute("90÷6=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "87÷5=17, 2", 2) | Out-Null
$d.Content.Find.Execute("62÷5=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "71÷4=17, 3", 2) | Out-Null
$d.Content.Find.Execute("42÷9=4, 6", $true, $false, $false, $false, $false, $true, 1, $false, "11÷8=1, 3", 2) | Out-Null
$d.Content.Find.Execute("75÷3=25, 0", $true, $false, $false, $false, $false, $true, 1, $false, "15÷8=1, 7", 2) | Out-Null
$d.Content.Find.Execute("66÷8=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "68÷6=11, 2", 2) | Out-Null
$d.Content.Find.Execute("22÷7=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "38÷4=9, 2", 2) | Out-Null
$d.Content.Find.Execute("22÷5=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "55÷4=13, 3", 2) | Out-Null
$d.Content.Find.Execute("81÷9=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "59÷3=19, 2", 2) | Out-Null
$d.Content.Find.Execute("36÷9=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "86÷5=17, 1", 2) | Out-Null
$d.Content.Find.Execute("78÷4=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "20÷9=2, 2", 2) | Out-Null
$d.Content.Find.Execute("74÷9=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "30÷7=4, 2", 2) | Out-Null
$d.Content.Find.Execute("91÷6=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "84÷9=9, 3", 2) | Out-Null
$d.Content.Find.Execute("53÷5=10, 3", $true, $false, $false, $false, $false, $true, 1, $false, "83÷8=10, 3", 2) | Out-Null
$d.Content.Find.Execute("96÷3=32, 0", $true, $false, $false, $false, $false, $true, 1, $false, "13÷5=2, 3", 2) | Out-Null
$d.Content.Find.Execute("49÷4=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "14÷6=2, 2", 2) | Out-Null
$d.Content.Find.Execute("99÷2=49, 1", $true, $false, $false, $false, $false, $true, 1, $false, "79÷3=26, 1", 2) | Out-Null
